$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value2 = '[ Nicolas%Hoertel%null%2,   Martin%Blachier%null%1,   Carlos%Blanco%null%1,   Mark%Olfson%null%1,   Marc%Massetti%null%1,   Marina Sánchez%Rico%null%1,   Frédéric%Limosin%null%1,   Henri%Leleu%null%1,  Nicolas%Hoertel%null%0,  Martin%Blachier%null%1,  Carlos%Blanco%null%1,  Mark%Olfson%null%1,  Marc%Massetti%null%1,  Marina Sánchez%Rico%null%1,  Frédéric%Limosin%null%1,  Henri%Leleu%null%1]'
$ws.Range("E5").Value2 = '[Tuo%Ji%NULL%1,  Hai-Lian%Chen%NULL%1,  Jing%Xu%NULL%1,  Ling-Ning%Wu%NULL%1,  Jie-Jia%Li%NULL%1,  Kai%Chen%NULL%1,  Gang%Qin%tonygqin@ntu.edu.cn%1]'
$ws.Range("I5").Value2 = '_PMC'
$ws.Range("E7").Value2 = '[Hien%Lau%hlau2@uci.edu%1,  Veria%Khosrawipour%veriakhosrawipour@yahoo.de%1,  Piotr%Kocbach%piotr.kocbach@uwm.edu.pl%1,  Agata%Mikolajczyk%agata.mikolajczyk@upwr.edu.pl%1,  Justyna%Schubert%justyna.schubert@upwr.edu.pl%1,  Jacek%Bania%jacek.bania@upwr.edu.pl%1,  Tanja%Khosrawipour%tkhosrawipour@gmail.com%1]'
$ws.Range("I7").Value2 = '_PMC'
$ws.Range("E8").Value2 = '[Carlo%Signorelli%NULL%1,  Thea%Scognamiglio%NULL%1,  Anna%Odone%NULL%1]'
$ws.Range("I8").Value2 = '_PMC'
$ws.Range("I10").Value2 = '_elsevier_PMC'
$ws.Range("E11").Value2 = '[Nadya%Johanna%NULL%1,  Henrico%Citrawijaya%NULL%1,  Grace%Wangge%NULL%1]'
$ws.Range("I11").Value2 = '_PMC'
$ws.Range("E13").Value2 = '[Biao%Tang%NULL%1,  Fan%Xia%NULL%1,  Sanyi%Tang%NULL%1,  Nicola Luigi%Bragazzi%NULL%1,  Qian%Li%NULL%1,  Xiaodan%Sun%NULL%1,  Juhua%Liang%NULL%1,  Yanni%Xiao%yxiao@mail.xjtu.edu.cn%1,  Jianhong%Wu%wujh@yorku.ca%1]'
$ws.Range("I13").Value2 = '_elsevier_PMC'
